$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 15.981972
$ws.Cells.Item(2, 8).Value = 47.945916
$ws.Cells.Item(2, 9).Value = 0.1372507760882863
$ws.Cells.Item(2, 10).Value = 0.1372507760882863
$ws.Cells.Item(2, 13).Value = 68.63737500000001
$ws.Cells.Item(2, 14).Value = 205.912125
$ws.Cells.Item(2, 15).Value = 0.5415701538216162
$ws.Cells.Item(2, 16).Value = 0.5415701538216162
$ws.Cells.Item(2, 17).Value = 1096.9606054035
$ws.Cells.Item(2, 18).Value = 9872.6454486315
$ws.Cells.Item(2, 19).Value = 0.0743309239182694
$ws.Cells.Item(2, 20).Value = 0.0743309239182694
$ws.Cells.Item(3, 7).Value = 15.981972
$ws.Cells.Item(3, 8).Value = 47.945916
$ws.Cells.Item(3, 9).Value = 0.1372507760882863
$ws.Cells.Item(3, 10).Value = 0.1372507760882863
$ws.Cells.Item(3, 15).Value = 0.08718851262838957
$ws.Cells.Item(3, 16).Value = 0.08718851262838957
$ws.Cells.Item(3, 17).Value = 176.601983920576
$ws.Cells.Item(3, 18).Value = 1589.417855285184
$ws.Cells.Item(3, 19).Value = 0.01196669102422982
$ws.Cells.Item(3, 20).Value = 0.01196669102422982
$ws.Cells.Item(4, 7).Value = 15.981972
$ws.Cells.Item(4, 8).Value = 47.945916
$ws.Cells.Item(4, 9).Value = 0.1372507760882863
$ws.Cells.Item(4, 10).Value = 0.1372507760882863
$ws.Cells.Item(4, 13).Value = 16.21089566666667
$ws.Cells.Item(4, 14).Value = 48.632687
$ws.Cells.Item(4, 15).Value = 0.1279089892319285
$ws.Cells.Item(4, 16).Value = 0.1279089892319285
$ws.Cells.Item(4, 17).Value = 259.082080639588
$ws.Cells.Item(4, 18).Value = 2331.738725756292
$ws.Cells.Item(4, 19).Value = 0.01755560804075044
$ws.Cells.Item(4, 20).Value = 0.01755560804075044
$ws.Cells.Item(5, 7).Value = 15.981972
$ws.Cells.Item(5, 8).Value = 47.945916
$ws.Cells.Item(5, 9).Value = 0.1372507760882863
$ws.Cells.Item(5, 10).Value = 0.1372507760882863
$ws.Cells.Item(5, 13).Value = 20.32546233333333
$ws.Cells.Item(5, 14).Value = 60.976387
$ws.Cells.Item(5, 15).Value = 0.1603741949973873
$ws.Cells.Item(5, 16).Value = 0.1603741949973873
$ws.Cells.Item(5, 17).Value = 324.840969898388
$ws.Cells.Item(5, 18).Value = 2923.568729085492
$ws.Cells.Item(5, 19).Value = 0.02201148272792557
$ws.Cells.Item(5, 20).Value = 0.02201148272792557
$ws.Cells.Item(6, 7).Value = 15.981972
$ws.Cells.Item(6, 8).Value = 47.945916
$ws.Cells.Item(6, 9).Value = 0.1372507760882863
$ws.Cells.Item(6, 10).Value = 0.1372507760882863
$ws.Cells.Item(6, 13).Value = 10.513928
$ws.Cells.Item(6, 14).Value = 31.541784
$ws.Cells.Item(6, 15).Value = 0.08295814932067838
$ws.Cells.Item(6, 16).Value = 0.08295814932067838
$ws.Cells.Item(6, 17).Value = 168.033302906016
$ws.Cells.Item(6, 18).Value = 1512.299726154144
$ws.Cells.Item(6, 19).Value = 0.01138607037711105
$ws.Cells.Item(6, 20).Value = 0.01138607037711105
$ws.Cells.Item(7, 8).Value = 69.213024
$ws.Cells.Item(7, 9).Value = 0.1981303529463737
$ws.Cells.Item(7, 10).Value = 0.1981303529463737
$ws.Cells.Item(7, 13).Value = 68.63737500000001
$ws.Cells.Item(7, 14).Value = 205.912125
$ws.Cells.Item(7, 15).Value = 0.5415701538216162
$ws.Cells.Item(7, 16).Value = 0.5415701538216162
$ws.Cells.Item(7, 17).Value = 1583.533427724
$ws.Cells.Item(7, 18).Value = 14251.800849516
$ws.Cells.Item(7, 19).Value = 0.1073014857218987
$ws.Cells.Item(7, 20).Value = 0.1073014857218987
$ws.Cells.Item(8, 8).Value = 69.213024
$ws.Cells.Item(8, 9).Value = 0.1981303529463737
$ws.Cells.Item(8, 10).Value = 0.1981303529463737
$ws.Cells.Item(8, 15).Value = 0.08718851262838957
$ws.Cells.Item(8, 16).Value = 0.08718851262838957
$ws.Cells.Item(8, 19).Value = 0.01727469077993218
$ws.Cells.Item(8, 20).Value = 0.01727469077993218
$ws.Cells.Item(9, 8).Value = 69.213024
$ws.Cells.Item(9, 9).Value = 0.1981303529463737
$ws.Cells.Item(9, 10).Value = 0.1981303529463737
$ws.Cells.Item(9, 13).Value = 16.21089566666667
$ws.Cells.Item(9, 14).Value = 48.632687
$ws.Cells.Item(9, 15).Value = 0.1279089892319285
$ws.Cells.Item(9, 16).Value = 0.1279089892319285
$ws.Cells.Item(9, 17).Value = 374.0017036128321
$ws.Cells.Item(9, 18).Value = 3366.015332515488
$ws.Cells.Item(9, 19).Value = 0.0253426531815359
$ws.Cells.Item(9, 20).Value = 0.0253426531815359
$ws.Cells.Item(10, 8).Value = 69.213024
$ws.Cells.Item(10, 9).Value = 0.1981303529463737
$ws.Cells.Item(10, 10).Value = 0.1981303529463737
$ws.Cells.Item(10, 13).Value = 20.32546233333333
$ws.Cells.Item(10, 14).Value = 60.976387
$ws.Cells.Item(10, 15).Value = 0.1603741949973873
$ws.Cells.Item(10, 16).Value = 0.1603741949973873
$ws.Cells.Item(10, 17).Value = 468.9289040960321
$ws.Cells.Item(10, 18).Value = 4220.360136864289
$ws.Cells.Item(10, 19).Value = 0.03177499585832291
$ws.Cells.Item(10, 20).Value = 0.03177499585832291
$ws.Cells.Item(11, 8).Value = 69.213024
$ws.Cells.Item(11, 9).Value = 0.1981303529463737
$ws.Cells.Item(11, 10).Value = 0.1981303529463737
$ws.Cells.Item(11, 13).Value = 10.513928
$ws.Cells.Item(11, 14).Value = 31.541784
$ws.Cells.Item(11, 15).Value = 0.08295814932067838
$ws.Cells.Item(11, 16).Value = 0.08295814932067838
$ws.Cells.Item(11, 17).Value = 242.566916999424
$ws.Cells.Item(11, 18).Value = 2183.102252994816
$ws.Cells.Item(11, 19).Value = 0.01643652740468398
$ws.Cells.Item(11, 20).Value = 0.01643652740468398
$ws.Cells.Item(12, 7).Value = 40.09539033333333
$ws.Cells.Item(12, 8).Value = 120.286171
$ws.Cells.Item(12, 9).Value = 0.3443331924754199
$ws.Cells.Item(12, 10).Value = 0.3443331924754199
$ws.Cells.Item(12, 13).Value = 68.63737500000001
$ws.Cells.Item(12, 14).Value = 205.912125
$ws.Cells.Item(12, 15).Value = 0.5415701538216162
$ws.Cells.Item(12, 16).Value = 0.5415701538216162
$ws.Cells.Item(12, 17).Value = 2752.042342080375
$ws.Cells.Item(12, 18).Value = 24768.38107872337
$ws.Cells.Item(12, 19).Value = 0.1864805800148013
$ws.Cells.Item(12, 20).Value = 0.1864805800148013
$ws.Cells.Item(13, 7).Value = 40.09539033333333
$ws.Cells.Item(13, 8).Value = 120.286171
$ws.Cells.Item(13, 9).Value = 0.3443331924754199
$ws.Cells.Item(13, 10).Value = 0.3443331924754199
$ws.Cells.Item(13, 15).Value = 0.08718851262838957
$ws.Cells.Item(13, 16).Value = 0.08718851262838957
$ws.Cells.Item(13, 17).Value = 443.0570569724783
$ws.Cells.Item(13, 18).Value = 3987.513512752304
$ws.Cells.Item(13, 19).Value = 0.03002189890051684
$ws.Cells.Item(13, 20).Value = 0.03002189890051684
$ws.Cells.Item(14, 7).Value = 40.09539033333333
$ws.Cells.Item(14, 8).Value = 120.286171
$ws.Cells.Item(14, 9).Value = 0.3443331924754199
$ws.Cells.Item(14, 10).Value = 0.3443331924754199
$ws.Cells.Item(14, 13).Value = 16.21089566666667
$ws.Cells.Item(14, 14).Value = 48.632687
$ws.Cells.Item(14, 15).Value = 0.1279089892319285
$ws.Cells.Item(14, 16).Value = 0.1279089892319285
$ws.Cells.Item(14, 17).Value = 649.982189407942
$ws.Cells.Item(14, 18).Value = 5849.839704671477
$ws.Cells.Item(14, 19).Value = 0.04404331060853404
$ws.Cells.Item(14, 20).Value = 0.04404331060853404
$ws.Cells.Item(15, 7).Value = 40.09539033333333
$ws.Cells.Item(15, 8).Value = 120.286171
$ws.Cells.Item(15, 9).Value = 0.3443331924754199
$ws.Cells.Item(15, 10).Value = 0.3443331924754199
$ws.Cells.Item(15, 13).Value = 20.32546233333333
$ws.Cells.Item(15, 14).Value = 60.976387
$ws.Cells.Item(15, 15).Value = 0.1603741949973873
$ws.Cells.Item(15, 16).Value = 0.1603741949973873
$ws.Cells.Item(15, 17).Value = 814.9573459604642
$ws.Cells.Item(15, 18).Value = 7334.616113644177
$ws.Cells.Item(15, 19).Value = 0.05522215855412589
$ws.Cells.Item(15, 20).Value = 0.05522215855412589
$ws.Cells.Item(16, 7).Value = 40.09539033333333
$ws.Cells.Item(16, 8).Value = 120.286171
$ws.Cells.Item(16, 9).Value = 0.3443331924754199
$ws.Cells.Item(16, 10).Value = 0.3443331924754199
$ws.Cells.Item(16, 13).Value = 10.513928
$ws.Cells.Item(16, 14).Value = 31.541784
$ws.Cells.Item(16, 15).Value = 0.08295814932067838
$ws.Cells.Item(16, 16).Value = 0.08295814932067838
$ws.Cells.Item(16, 17).Value = 421.5600470965627
$ws.Cells.Item(16, 18).Value = 3794.040423869064
$ws.Cells.Item(16, 19).Value = 0.02856524439744177
$ws.Cells.Item(16, 20).Value = 0.02856524439744177
$ws.Cells.Item(17, 7).Value = 8.831340666666666
$ws.Cells.Item(17, 8).Value = 26.494022
$ws.Cells.Item(17, 9).Value = 0.07584222775512579
$ws.Cells.Item(17, 10).Value = 0.07584222775512579
$ws.Cells.Item(17, 13).Value = 68.63737500000001
$ws.Cells.Item(17, 14).Value = 205.912125
$ws.Cells.Item(17, 15).Value = 0.5415701538216162
$ws.Cells.Item(17, 16).Value = 0.5415701538216162
$ws.Cells.Item(17, 17).Value = 606.16004109075
$ws.Cells.Item(17, 18).Value = 5455.44036981675
$ws.Cells.Item(17, 19).Value = 0.04107388695151753
$ws.Cells.Item(17, 20).Value = 0.04107388695151753
$ws.Cells.Item(18, 7).Value = 8.831340666666666
$ws.Cells.Item(18, 8).Value = 26.494022
$ws.Cells.Item(18, 9).Value = 0.07584222775512579
$ws.Cells.Item(18, 10).Value = 0.07584222775512579
$ws.Cells.Item(18, 15).Value = 0.08718851262838957
$ws.Cells.Item(18, 16).Value = 0.08718851262838957
$ws.Cells.Item(18, 17).Value = 97.58697377343644
$ws.Cells.Item(18, 18).Value = 878.2827639609279
$ws.Cells.Item(18, 19).Value = 0.006612571032392984
$ws.Cells.Item(18, 20).Value = 0.006612571032392984
$ws.Cells.Item(19, 7).Value = 8.831340666666666
$ws.Cells.Item(19, 8).Value = 26.494022
$ws.Cells.Item(19, 9).Value = 0.07584222775512579
$ws.Cells.Item(19, 10).Value = 0.07584222775512579
$ws.Cells.Item(19, 13).Value = 16.21089566666667
$ws.Cells.Item(19, 14).Value = 48.632687
$ws.Cells.Item(19, 15).Value = 0.1279089892319285
$ws.Cells.Item(19, 16).Value = 0.1279089892319285
$ws.Cells.Item(19, 17).Value = 143.1639421441238
$ws.Cells.Item(19, 18).Value = 1288.475479297114
$ws.Cells.Item(19, 19).Value = 0.009700902693255855
$ws.Cells.Item(19, 20).Value = 0.009700902693255855
$ws.Cells.Item(20, 7).Value = 8.831340666666666
$ws.Cells.Item(20, 8).Value = 26.494022
$ws.Cells.Item(20, 9).Value = 0.07584222775512579
$ws.Cells.Item(20, 10).Value = 0.07584222775512579
$ws.Cells.Item(20, 13).Value = 20.32546233333333
$ws.Cells.Item(20, 14).Value = 60.976387
$ws.Cells.Item(20, 15).Value = 0.1603741949973873
$ws.Cells.Item(20, 16).Value = 0.1603741949973873
$ws.Cells.Item(20, 17).Value = 179.5010820731682
$ws.Cells.Item(20, 18).Value = 1615.509738658514
$ws.Cells.Item(20, 19).Value = 0.01216313622303681
$ws.Cells.Item(20, 20).Value = 0.01216313622303681
$ws.Cells.Item(21, 7).Value = 8.831340666666666
$ws.Cells.Item(21, 8).Value = 26.494022
$ws.Cells.Item(21, 9).Value = 0.07584222775512579
$ws.Cells.Item(21, 10).Value = 0.07584222775512579
$ws.Cells.Item(21, 13).Value = 10.513928
$ws.Cells.Item(21, 14).Value = 31.541784
$ws.Cells.Item(21, 15).Value = 0.08295814932067838
$ws.Cells.Item(21, 16).Value = 0.08295814932067838
$ws.Cells.Item(21, 17).Value = 92.85207991280532
$ws.Cells.Item(21, 18).Value = 835.668719215248
$ws.Cells.Item(21, 19).Value = 0.006291730854922624
$ws.Cells.Item(21, 20).Value = 0.006291730854922624
$ws.Cells.Item(22, 7).Value = 28.463871
$ws.Cells.Item(22, 8).Value = 85.39161300000001
$ws.Cells.Item(22, 9).Value = 0.2444434507347945
$ws.Cells.Item(22, 10).Value = 0.2444434507347945
$ws.Cells.Item(22, 13).Value = 68.63737500000001
$ws.Cells.Item(22, 14).Value = 205.912125
$ws.Cells.Item(22, 15).Value = 0.5415701538216162
$ws.Cells.Item(22, 16).Value = 0.5415701538216162
$ws.Cells.Item(22, 17).Value = 1953.685387778625
$ws.Cells.Item(22, 18).Value = 17583.16849000763
$ws.Cells.Item(22, 19).Value = 0.1323832772151293
$ws.Cells.Item(22, 20).Value = 0.1323832772151293
$ws.Cells.Item(23, 7).Value = 28.463871
$ws.Cells.Item(23, 8).Value = 85.39161300000001
$ws.Cells.Item(23, 9).Value = 0.2444434507347945
$ws.Cells.Item(23, 10).Value = 0.2444434507347945
$ws.Cells.Item(23, 15).Value = 0.08718851262838957
$ws.Cells.Item(23, 16).Value = 0.08718851262838957
$ws.Cells.Item(23, 17).Value = 314.527899852368
$ws.Cells.Item(23, 18).Value = 2830.751098671312
$ws.Cells.Item(23, 19).Value = 0.02131266089131775
$ws.Cells.Item(23, 20).Value = 0.02131266089131776
$ws.Cells.Item(24, 7).Value = 28.463871
$ws.Cells.Item(24, 8).Value = 85.39161300000001
$ws.Cells.Item(24, 9).Value = 0.2444434507347945
$ws.Cells.Item(24, 10).Value = 0.2444434507347945
$ws.Cells.Item(24, 13).Value = 16.21089566666667
$ws.Cells.Item(24, 14).Value = 48.632687
$ws.Cells.Item(24, 15).Value = 0.1279089892319285
$ws.Cells.Item(24, 16).Value = 0.1279089892319285
$ws.Cells.Item(24, 17).Value = 461.4248430504591
$ws.Cells.Item(24, 18).Value = 4152.823587454131
$ws.Cells.Item(24, 19).Value = 0.03126651470785227
$ws.Cells.Item(24, 20).Value = 0.03126651470785228
$ws.Cells.Item(25, 7).Value = 28.463871
$ws.Cells.Item(25, 8).Value = 85.39161300000001
$ws.Cells.Item(25, 9).Value = 0.2444434507347945
$ws.Cells.Item(25, 10).Value = 0.2444434507347945
$ws.Cells.Item(25, 13).Value = 20.32546233333333
$ws.Cells.Item(25, 14).Value = 60.976387
$ws.Cells.Item(25, 15).Value = 0.1603741949973873
$ws.Cells.Item(25, 16).Value = 0.1603741949973873
$ws.Cells.Item(25, 17).Value = 578.541337871359
$ws.Cells.Item(25, 18).Value = 5206.872040842231
$ws.Cells.Item(25, 19).Value = 0.03920242163397618
$ws.Cells.Item(25, 20).Value = 0.03920242163397618
$ws.Cells.Item(26, 7).Value = 28.463871
$ws.Cells.Item(26, 8).Value = 85.39161300000001
$ws.Cells.Item(26, 9).Value = 0.2444434507347945
$ws.Cells.Item(26, 10).Value = 0.2444434507347945
$ws.Cells.Item(26, 13).Value = 10.513928
$ws.Cells.Item(26, 14).Value = 31.541784
$ws.Cells.Item(26, 15).Value = 0.08295814932067838
$ws.Cells.Item(26, 16).Value = 0.08295814932067838
$ws.Cells.Item(26, 17).Value = 299.267090295288
$ws.Cells.Item(26, 18).Value = 2693.403812657592
$ws.Cells.Item(26, 19).Value = 0.02027857628651897
$ws.Cells.Item(26, 20).Value = 0.02027857628651897
